$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new localization row for the "full screen" setting option.
$ws.Range("A19").Value = "full_screen"
$ws.Range("B19").Value = "全屏"
$ws.Range("C19").Value = "FULL SCREEN"
$ws.Range("D19").Value = "全屏"

# Match the author's cursor position when the file was saved.
$ws.Range("E16").Select()
